$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.251.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.84%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.821.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.02%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4470"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3765"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07394"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8794"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.85"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.52%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.822.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.708"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.423"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07065"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.00%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008805"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.26%  "

$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.251.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.344"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.958"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.286"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.348"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08868"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.7915"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.87%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.194"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.568"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.925"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.21%  "

$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.107"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01975"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05265"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.290"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5280"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.869"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.80%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.329"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +18.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1701"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.636"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5042"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.686"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.17%  "

$ws.Range("E49").Value = "  +0.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06382"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "66.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.58%  "
